$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.894.91'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '3.509.85'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.54'
$ws.Range('E5').Value = '  +3.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.27'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.216'
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.664'
$ws.Range('E10').Value = '  +2.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.64'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000309'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.62'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').Value = '4.072.66'
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '617.08'
$ws.Range('E15').Value = '  +8.18%  '
$ws.Range('D16').Value = '70.000.04'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.74'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.92'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').Value = '3.509.99'
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.53'
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '105.91'
$ws.Range('E23').Value = '  +12.40%  '
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.01'
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.05'
$ws.Range('E26').Value = '  +4.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.00'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.88'
$ws.Range('E28').Value = '  +5.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.24'
$ws.Range('E29').Value = '  +5.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.50'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('B32').Value = 'dogwifhat'
$ws.Range('C32').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.09'
$ws.Range('E32').Value = '  +5.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.97'
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '3.709.67'
$ws.Range('E35').Value = '  +2.33%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.10'
$ws.Range('E36').Value = '  -4.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '516.38'
$ws.Range('E38').Value = '  -1.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.392'
$ws.Range('E39').Value = '  -4.24%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0784'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.59'
$ws.Range('E41').Value = '  +1.70%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.78'
$ws.Range('E42').Value = '  -4.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.138'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0464'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.88'
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.141'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').Value = '  -4.17%  '
$ws.Range('E48').Value = '  -4.94%  '
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.11'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('E51').Value = '  -6.29%  '
